$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.258.92"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.656.27"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5243"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06367"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.602"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.667.37"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "1.886.16"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5643"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "0.0₅8263"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "26.265.31"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "192.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1206"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.285"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.513"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05665"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.281"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.512"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.587"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.12%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9521"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5778"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.986"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.566"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8466"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").Value = "1.015.78"
$ws.Range("E45").Value = "  -5.72%  "
$ws.Range("D46").Value = "1.796.46"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05353"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.028"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4351"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "

Write-Output "Applied 90 cell updates"
